$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-09-05 19:13:37"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-09-05 19:13:33"
$wsZhCn.Range("K2").Value = "2016-09-05 19:13:57"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-09-05 19:13:37"
$wsDeDe.Range("K2").Value = "2016-09-05 19:14:11"
